# Append four new daily rows (71-74) to the "100 Error Counts" sheet,
# matching the data that was uploaded in the commit, and move the
# viewport/selection down to show the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows: row, Date (serial), Total Count, Session Timeout Errors, Errors Requiring Analysis
$newRows = @(
    @(71, 46031, 614, 597, 17),
    @(72, 46034, 748, 722, 26),
    @(73, 46035, 793, 744, 49),
    @(74, 46036, 637, 613, 24)
)

foreach ($row in $newRows) {
    $r = $row[0]

    # Column A carries the date number format ("m/d/yy" resolves to the
    # same built-in short-date format already used by the rows above it).
    $ws.Cells.Item($r, 1).NumberFormat = "m/d/yy"
    $ws.Cells.Item($r, 1).Value = $row[1]

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Scroll the view down and select the newly added block, same as the
# author leaving the selection on the last entered rows.
$ws.Range("A69:D74").Select()
$excel.ActiveWindow.ScrollRow = 59
